$p = $ppt.ActivePresentation

# Remove the last two slides: "Panel design" (slide 4) and "Panel assembly" (slide 5)
$p.Slides.Item(5).Delete()
$p.Slides.Item(4).Delete()
